## Adds the new "Co" (AG) and "Yb" (AH) trace-element columns to the
## Total-Quebradagrande sheet, fills in the available analyses for rows
## 11-30, and nudges the selection to the last edited cell -- mirroring
## the "adding cleaned .gas project for ioGas and cleaned tables with Th
## and Co (except for Nivia, 2006)" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -----------------------------------------------
$ws.Range("AG1").Value = "Co"
$ws.Range("AH1").Value = "Yb"

# --- New Co / Yb analyses, rows 11-16 ----------------------------------
# These rows keep the same "0" / "0.00" number formats already used by
# the rest of that block (columns C..AF for those rows).
$rows1116 = @(
  @{R=11; Co=13;   Yb=1.77},
  @{R=12; Co=15;   Yb=1.25},
  @{R=13; Co=40;   Yb=3.72},
  @{R=14; Co=14;   Yb=1.35},
  @{R=15; Co=30;   Yb=1.43},
  @{R=16; Co=19;   Yb=1.74}
)
foreach ($item in $rows1116) {
  $r = $item.R
  $agCell = $ws.Range("AG$r")
  $agCell.Value = $item.Co
  $agCell.NumberFormat = "0"

  $ahCell = $ws.Range("AH$r")
  $ahCell.Value = $item.Yb
  $ahCell.NumberFormat = "0.00"
}

# Rows 17-30 use the plain general format, matching the rest of that block.
$rows1730 = @(
  @{R=17; Co=57.7;               Yb=2.2999999999999998},
  @{R=18; Co=35.1;               Yb=2.4},
  @{R=19; Co=59.3;               Yb=3.7},
  @{R=20; Co=52.6;               Yb=4},
  @{R=21; Co=59.3;               Yb=3.3},
  @{R=22; Co=55.2;               Yb=3.7},
  @{R=23; Co=62.5;               Yb=4.2},
  @{R=24; Co=52;                 Yb=3.2},
  @{R=25; Co=56.9;               Yb=3.1},
  @{R=26; Co=58;                 Yb=3},
  @{R=27; Co=62;                 Yb=2.6},
  @{R=28; Co=45.1;               Yb=3.1},
  @{R=29; Co=40.799999999999997; Yb=1.5},
  @{R=30; Co=42.4;               Yb=1.8}
)
foreach ($item in $rows1730) {
  $r = $item.R
  $ws.Range("AG$r").Value = $item.Co
  $ws.Range("AH$r").Value = $item.Yb
}

# --- Drop the underline styling from the two "Rod..." sample labels ----
# (B10 and A29 were underlined before; the cleaned table removes it while
# keeping B15's plain style and A29's "Rod"@ custom number format.)
$ws.Range("B10").Font.Underline = $false
$ws.Range("A29").Font.Underline = $false

# --- Selection moves to the last cell touched, like the source file ----
$ws.Range("AH15").Select() | Out-Null
